$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Merge the two split "-----" runs (Week 2 -> Week 3 separator)
#    into a single run of 51 dashes.
# ---------------------------------------------------------------
$oldDashes = "-----------------------------------------" + "----------"
$newDashes = "---------------------------------------------------"

$sepPara = $d.Paragraphs.Item(36)
$sepRange = $sepPara.Range
$sepRange.Find.ClearFormatting()
$sepRange.Find.Execute($oldDashes, $false, $false, $false, $false, $false, $true, 1, $false, $newDashes, 2)

# ---------------------------------------------------------------
# 2) Week 3 table: check both "Checked" cells ("x").
#    Row 2 cell needs only the "x" text (no bookmark involved).
# ---------------------------------------------------------------
$week3Table = $d.Tables.Item(3)
$week3Table.Cell(2, 2).Range.Text = "x"

# Row 3 cell needs the "x" text *and* the "_GoBack" bookmark that is
# being moved here from the Week 4 table below. Rebuild the whole
# paragraph (preserving its pPr / rsid markers) via InsertXML so the
# bookmark pair lands inside it.
$week3Para = $week3Table.Cell(3, 2).Range.Paragraphs.Item(1)
$week3ParaRange = $week3Para.Range
$week3Xml = '<?xml version="1.0"?><?mso-application progid="Word.Document"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p w:rsidR="008A032F" w:rsidRDefault="008A032F">' + `
  '<w:pPr><w:widowControl w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>' + `
  '<w:r><w:t>x</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
  '<w:bookmarkEnd w:id="0"/>' + `
  '</w:p>' + `
  '</w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'
$week3ParaRange.InsertXML($week3Xml)

# ---------------------------------------------------------------
# 3) Week 4 table: the second "Checked" cell used to own the
#    "_GoBack" bookmark; it has now moved to the Week 3 table above,
#    so remove it from here (cell stays textless).
# ---------------------------------------------------------------
$week4Table = $d.Tables.Item(4)
$week4Para = $week4Table.Cell(3, 2).Range.Paragraphs.Item(1)
$week4ParaRange = $week4Para.Range
$week4Xml = '<?xml version="1.0"?><?mso-application progid="Word.Document"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p w:rsidR="008A032F" w:rsidRDefault="008A032F">' + `
  '<w:pPr><w:widowControl w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>' + `
  '</w:p>' + `
  '</w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'
$week4ParaRange.InsertXML($week4Xml)

Write-Output "edit complete"
